$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update values in A2 and A3
$ws.Range("A2").Value = 71472
$ws.Range("A3").Value = 71475

# Remove row 4 entirely (was A4 = 71409)
$ws.Range("A4").EntireRow.Delete()

# Update the active selection to T4
$ws.Range("T4").Select()
